# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains an extra (blank) column between the
# existing "Late" and "Outstanding" columns (i.e. a new blank column is
# inserted immediately before the old "Late" column, pushing
# Late -> one column right, and Outstanding -> two columns right).
# The sheet also becomes the active/selected sheet & tab, with the
# selection left on R8 (instead of the "Input" sheet being active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"), shifting
# N -> O ("Late"), O -> P (blank spacer), P -> Q ("Outstanding").
$ws.Columns("N:N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab (this also clears
# tabSelected on the previously active "Input" sheet), and leave the
# selection on R8 to match the saved view state.
$ws.Select() | Out-Null
$ws.Range("R8").Select() | Out-Null
